$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = -11.62650000000001
$ws.Range("C7").Value = -11.6674
$ws.Range("B8").Value = 4.695200000000002
$ws.Range("A12").Value = -22.7526
$ws.Range("B12").Value = 6.043800000000001
$ws.Range("B14").Value = 8.599700000000006
$ws.Range("C19").Value = -13.37569999999999
$ws.Range("C21").Value = -13.1029
$ws.Range("B22").Value = 4.681800000000005
$ws.Range("C24").Value = -11.3859
